$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 27.73790633333333
$ws.Cells.Item(2, 8).Value = 83.213719
$ws.Cells.Item(2, 9).Value = 0.005442473085408622
$ws.Cells.Item(2, 10).Value = 0.005456294496964538
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 2.476839
$ws.Cells.Item(2, 14).Value = 7.430517
$ws.Cells.Item(2, 15).Value = 0.0525987134655237
$ws.Cells.Item(2, 16).Value = 0.05675564862155354
$ws.Cells.Item(2, 17).Value = 68.70232818474699
$ws.Cells.Item(2, 18).Value = 618.320953662723
$ws.Cells.Item(2, 19).Value = 0.0002862670823632328
$ws.Cells.Item(2, 20).Value = 0.0003096755332454355

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 27.73790633333333
$ws.Cells.Item(3, 8).Value = 83.213719
$ws.Cells.Item(3, 9).Value = 0.005442473085408622
$ws.Cells.Item(3, 10).Value = 0.005456294496964538
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 7.339638666666666
$ws.Cells.Item(3, 14).Value = 22.018916
$ws.Cells.Item(3, 15).Value = 0.1558662275458673
$ws.Cells.Item(3, 16).Value = 0.1681845098427879
$ws.Cells.Item(3, 17).Value = 203.5862098565115
$ws.Cells.Item(3, 18).Value = 1832.275888708604
$ws.Cells.Item(3, 19).Value = 0.0008482977483425585
$ws.Cells.Item(3, 20).Value = 0.0009176642155298819

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 27.73790633333333
$ws.Cells.Item(4, 8).Value = 83.213719
$ws.Cells.Item(4, 9).Value = 0.005442473085408622
$ws.Cells.Item(4, 10).Value = 0.005456294496964538
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 9.137454333333332
$ws.Cells.Item(4, 14).Value = 27.412363
$ws.Cells.Item(4, 15).Value = 0.1940450478546679
$ws.Cells.Item(4, 16).Value = 0.2093806450230146
$ws.Cells.Item(4, 17).Value = 253.4538524231107
$ws.Cells.Item(4, 18).Value = 2281.084671807997
$ws.Cells.Item(4, 19).Value = 0.001056084950305858
$ws.Cells.Item(4, 20).Value = 0.00114244246120996

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 27.73790633333333
$ws.Cells.Item(5, 8).Value = 83.213719
$ws.Cells.Item(5, 9).Value = 0.005442473085408622
$ws.Cells.Item(5, 10).Value = 0.005456294496964538
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 17.78856566666667
$ws.Cells.Item(5, 14).Value = 53.365697
$ws.Cells.Item(5, 15).Value = 0.3777620057111716
$ws.Cells.Item(5, 16).Value = 0.4076169595435007
$ws.Cells.Item(5, 17).Value = 493.4175682663492
$ws.Cells.Item(5, 18).Value = 4440.758114397143
$ws.Cells.Item(5, 19).Value = 0.002055959548773029
$ws.Cells.Item(5, 20).Value = 0.002224078173226619

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 27.73790633333333
$ws.Cells.Item(6, 8).Value = 83.213719
$ws.Cells.Item(6, 9).Value = 0.005442473085408622
$ws.Cells.Item(6, 10).Value = 0.005456294496964538
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 13).Value = 10.346848
$ws.Cells.Item(6, 14).Value = 20.693696
$ws.Cells.Item(6, 15).Value = 0.2197280054227695
$ws.Cells.Item(6, 16).Value = 0.1580622369691433
$ws.Cells.Item(6, 17).Value = 286.9999006692373
$ws.Cells.Item(6, 18).Value = 1721.999404015424
$ws.Cells.Item(6, 19).Value = 0.001195863755623942
$ws.Cells.Item(6, 20).Value = 0.0008624341137526413

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 45.15892033333333
$ws.Cells.Item(7, 8).Value = 135.476761
$ws.Cells.Item(7, 9).Value = 0.008860661851212738
$ws.Cells.Item(7, 10).Value = 0.008883163910879647
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 2.476839
$ws.Cells.Item(7, 14).Value = 7.430517
$ws.Cells.Item(7, 15).Value = 0.0525987134655237
$ws.Cells.Item(7, 16).Value = 0.05675564862155354
$ws.Cells.Item(7, 17).Value = 111.851375079493
$ws.Cells.Item(7, 18).Value = 1006.662375715437
$ws.Cells.Item(7, 19).Value = 0.0004660594138268356
$ws.Cells.Item(7, 20).Value = 0.0005041697295735506

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 45.15892033333333
$ws.Cells.Item(8, 8).Value = 135.476761
$ws.Cells.Item(8, 9).Value = 0.008860661851212738
$ws.Cells.Item(8, 10).Value = 0.008883163910879647
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 7.339638666666666
$ws.Cells.Item(8, 14).Value = 22.018916
$ws.Cells.Item(8, 15).Value = 0.1558662275458673
$ws.Cells.Item(8, 16).Value = 0.1681845098427879
$ws.Cells.Item(8, 17).Value = 331.4501578234529
$ws.Cells.Item(8, 18).Value = 2983.051420411076
$ws.Cells.Item(8, 19).Value = 0.00138107793630811
$ws.Cells.Item(8, 20).Value = 0.001494010568204437

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 45.15892033333333
$ws.Cells.Item(9, 8).Value = 135.476761
$ws.Cells.Item(9, 9).Value = 0.008860661851212738
$ws.Cells.Item(9, 10).Value = 0.008883163910879647
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 9.137454333333332
$ws.Cells.Item(9, 14).Value = 27.412363
$ws.Cells.Item(9, 15).Value = 0.1940450478546679
$ws.Cells.Item(9, 16).Value = 0.2093806450230146
$ws.Cells.Item(9, 17).Value = 412.6375722884714
$ws.Cells.Item(9, 18).Value = 3713.738150596243
$ws.Cells.Item(9, 19).Value = 0.001719367552942606
$ws.Cells.Item(9, 20).Value = 0.001859962589505145

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 45.15892033333333
$ws.Cells.Item(10, 8).Value = 135.476761
$ws.Cells.Item(10, 9).Value = 0.008860661851212738
$ws.Cells.Item(10, 10).Value = 0.008883163910879647
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 17.78856566666667
$ws.Cells.Item(10, 14).Value = 53.365697
$ws.Cells.Item(10, 15).Value = 0.3777620057111716
$ws.Cells.Item(10, 16).Value = 0.4076169595435007
$ws.Cells.Item(10, 17).Value = 803.3124197852686
$ws.Cells.Item(10, 18).Value = 7229.811778067418
$ws.Cells.Item(10, 19).Value = 0.003347221392842587
$ws.Cells.Item(10, 20).Value = 0.003620928264479315

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 45.15892033333333
$ws.Cells.Item(11, 8).Value = 135.476761
$ws.Cells.Item(11, 9).Value = 0.008860661851212738
$ws.Cells.Item(11, 10).Value = 0.008883163910879647
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 13).Value = 10.346848
$ws.Cells.Item(11, 14).Value = 20.693696
$ws.Cells.Item(11, 15).Value = 0.2197280054227695
$ws.Cells.Item(11, 16).Value = 0.1580622369691433
$ws.Cells.Item(11, 17).Value = 467.2524845331093
$ws.Cells.Item(11, 18).Value = 2803.514907198656
$ws.Cells.Item(11, 19).Value = 0.001946935555292599
$ws.Cells.Item(11, 20).Value = 0.001404092759117201

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 2612.668416333333
$ws.Cells.Item(12, 8).Value = 7838.005249
$ws.Cells.Item(12, 9).Value = 0.5126334109760676
$ws.Cells.Item(12, 10).Value = 0.5139352671798969
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 2.476839
$ws.Cells.Item(12, 14).Value = 7.430517
$ws.Cells.Item(12, 15).Value = 0.0525987134655237
$ws.Cells.Item(12, 16).Value = 0.05675564862155354
$ws.Cells.Item(12, 17).Value = 6471.159027642637
$ws.Cells.Item(12, 18).Value = 58240.43124878373
$ws.Cells.Item(12, 19).Value = 0.02696385789678423
$ws.Cells.Item(12, 20).Value = 0.02916872943828646

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 2612.668416333333
$ws.Cells.Item(13, 8).Value = 7838.005249
$ws.Cells.Item(13, 9).Value = 0.5126334109760676
$ws.Cells.Item(13, 10).Value = 0.5139352671798969
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 7.339638666666666
$ws.Cells.Item(13, 14).Value = 22.018916
$ws.Cells.Item(13, 15).Value = 0.1558662275458673
$ws.Cells.Item(13, 16).Value = 0.1681845098427879
$ws.Cells.Item(13, 17).Value = 19176.0421316989
$ws.Cells.Item(13, 18).Value = 172584.37918529
$ws.Cells.Item(13, 19).Value = 0.07990223588280984
$ws.Cells.Item(13, 20).Value = 0.08643595100157321

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 2612.668416333333
$ws.Cells.Item(14, 8).Value = 7838.005249
$ws.Cells.Item(14, 9).Value = 0.5126334109760676
$ws.Cells.Item(14, 10).Value = 0.5139352671798969
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 9.137454333333332
$ws.Cells.Item(14, 14).Value = 27.412363
$ws.Cells.Item(14, 15).Value = 0.1940450478546679
$ws.Cells.Item(14, 16).Value = 0.2093806450230146
$ws.Cells.Item(14, 17).Value = 23873.13834238815
$ws.Cells.Item(14, 18).Value = 214858.2450814934
$ws.Cells.Item(14, 19).Value = 0.09947397476475268
$ws.Cells.Item(14, 20).Value = 0.1076080977422021

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 2612.668416333333
$ws.Cells.Item(15, 8).Value = 7838.005249
$ws.Cells.Item(15, 9).Value = 0.5126334109760676
$ws.Cells.Item(15, 10).Value = 0.5139352671798969
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 17.78856566666667
$ws.Cells.Item(15, 14).Value = 53.365697
$ws.Cells.Item(15, 15).Value = 0.3777620057111716
$ws.Cells.Item(15, 16).Value = 0.4076169595435007
$ws.Cells.Item(15, 17).Value = 46475.62368917151
$ws.Cells.Item(15, 18).Value = 418280.6132025436
$ws.Cells.Item(15, 19).Value = 0.1936534255248786
$ws.Cells.Item(15, 20).Value = 0.2094887310100463

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 2612.668416333333
$ws.Cells.Item(16, 8).Value = 7838.005249
$ws.Cells.Item(16, 9).Value = 0.5126334109760676
$ws.Cells.Item(16, 10).Value = 0.5139352671798969
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 13).Value = 10.346848
$ws.Cells.Item(16, 14).Value = 20.693696
$ws.Cells.Item(16, 15).Value = 0.2197280054227695
$ws.Cells.Item(16, 16).Value = 0.1580622369691433
$ws.Cells.Item(16, 17).Value = 27032.88297820172
$ws.Cells.Item(16, 18).Value = 162197.2978692103
$ws.Cells.Item(16, 19).Value = 0.1126399169068422
$ws.Cells.Item(16, 20).Value = 0.08123375798778885

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 2372.267130666667
$ws.Cells.Item(17, 8).Value = 7116.801392
$ws.Cells.Item(17, 9).Value = 0.4654641145188886
$ws.Cells.Item(17, 10).Value = 0.4666461821176285
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 2.476839
$ws.Cells.Item(17, 14).Value = 7.430517
$ws.Cells.Item(17, 15).Value = 0.0525987134655237
$ws.Cells.Item(17, 16).Value = 0.05675564862155354
$ws.Cells.Item(17, 17).Value = 5875.723747653296
$ws.Cells.Item(17, 18).Value = 52881.51372887967
$ws.Cells.Item(17, 19).Value = 0.02448281358806273
$ws.Cells.Item(17, 20).Value = 0.0264848067428576

$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 2372.267130666667
$ws.Cells.Item(18, 8).Value = 7116.801392
$ws.Cells.Item(18, 9).Value = 0.4654641145188886
$ws.Cells.Item(18, 10).Value = 0.4666461821176285
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 7.339638666666666
$ws.Cells.Item(18, 14).Value = 22.018916
$ws.Cells.Item(18, 15).Value = 0.1558662275458673
$ws.Cells.Item(18, 16).Value = 0.1681845098427879
$ws.Cells.Item(18, 17).Value = 17411.58355990345
$ws.Cells.Item(18, 18).Value = 156704.2520391311
$ws.Cells.Item(18, 19).Value = 0.07255013558803672
$ws.Cells.Item(18, 20).Value = 0.07848265940946171

$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 2372.267130666667
$ws.Cells.Item(19, 8).Value = 7116.801392
$ws.Cells.Item(19, 9).Value = 0.4654641145188886
$ws.Cells.Item(19, 10).Value = 0.4666461821176285
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 13).Value = 9.137454333333332
$ws.Cells.Item(19, 14).Value = 27.412363
$ws.Cells.Item(19, 15).Value = 0.1940450478546679
$ws.Cells.Item(19, 16).Value = 0.2093806450230146
$ws.Cells.Item(19, 17).Value = 21676.48257293436
$ws.Cells.Item(19, 18).Value = 195088.3431564093
$ws.Cells.Item(19, 19).Value = 0.09032100637644837
$ws.Cells.Item(19, 20).Value = 0.09770667860931619

$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 2372.267130666667
$ws.Cells.Item(20, 8).Value = 7116.801392
$ws.Cells.Item(20, 9).Value = 0.4654641145188886
$ws.Cells.Item(20, 10).Value = 0.4666461821176285
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 13).Value = 17.78856566666667
$ws.Cells.Item(20, 14).Value = 53.365697
$ws.Cells.Item(20, 15).Value = 0.3777620057111716
$ws.Cells.Item(20, 16).Value = 0.4076169595435007
$ws.Cells.Item(20, 17).Value = 42199.22963273891
$ws.Cells.Item(20, 18).Value = 379793.0666946503
$ws.Cells.Item(20, 19).Value = 0.1758346574872298
$ws.Cells.Item(20, 20).Value = 0.1902128979373705

$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 2372.267130666667
$ws.Cells.Item(21, 8).Value = 7116.801392
$ws.Cells.Item(21, 9).Value = 0.4654641145188886
$ws.Cells.Item(21, 10).Value = 0.4666461821176285
$ws.Cells.Item(21, 11).Value = 2
$ws.Cells.Item(21, 13).Value = 10.346848
$ws.Cells.Item(21, 14).Value = 20.693696
$ws.Cells.Item(21, 15).Value = 0.2197280054227695
$ws.Cells.Item(21, 16).Value = 0.1580622369691433
$ws.Cells.Item(21, 17).Value = 24545.48741640414
$ws.Cells.Item(21, 18).Value = 147272.9244984248
$ws.Cells.Item(21, 19).Value = 0.102275501479111
$ws.Cells.Item(21, 20).Value = 0.07375913941862261

$ws.Cells.Item(22, 5).Value = 2
$ws.Cells.Item(22, 7).Value = 38.730512
$ws.Cells.Item(22, 8).Value = 77.461024
$ws.Cells.Item(22, 9).Value = 0.00759933956842245
$ws.Cells.Item(22, 10).Value = 0.005079092294630384
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 2.476839
$ws.Cells.Item(22, 14).Value = 7.430517
$ws.Cells.Item(22, 15).Value = 0.0525987134655237
$ws.Cells.Item(22, 16).Value = 0.05675564862155354
$ws.Cells.Item(22, 17).Value = 95.929242611568
$ws.Cells.Item(22, 18).Value = 575.5754556694079
$ws.Cells.Item(22, 19).Value = 0.000399715484486669
$ws.Cells.Item(22, 20).Value = 0.0002882671775904821

$ws.Cells.Item(23, 5).Value = 2
$ws.Cells.Item(23, 7).Value = 38.730512
$ws.Cells.Item(23, 8).Value = 77.461024
$ws.Cells.Item(23, 9).Value = 0.00759933956842245
$ws.Cells.Item(23, 10).Value = 0.005079092294630384
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 7.339638666666666
$ws.Cells.Item(23, 14).Value = 22.018916
$ws.Cells.Item(23, 15).Value = 0.1558662275458673
$ws.Cells.Item(23, 16).Value = 0.1681845098427879
$ws.Cells.Item(23, 17).Value = 284.2679634549973
$ws.Cells.Item(23, 18).Value = 1705.607780729984
$ws.Cells.Item(23, 19).Value = 0.001184480390370046
$ws.Cells.Item(23, 20).Value = 0.0008542246480186921

$ws.Cells.Item(24, 5).Value = 2
$ws.Cells.Item(24, 7).Value = 38.730512
$ws.Cells.Item(24, 8).Value = 77.461024
$ws.Cells.Item(24, 9).Value = 0.00759933956842245
$ws.Cells.Item(24, 10).Value = 0.005079092294630384
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 9.137454333333332
$ws.Cells.Item(24, 14).Value = 27.412363
$ws.Cells.Item(24, 15).Value = 0.1940450478546679
$ws.Cells.Item(24, 16).Value = 0.2093806450230146
$ws.Cells.Item(24, 17).Value = 353.8982847066186
$ws.Cells.Item(24, 18).Value = 2123.389708239712
$ws.Cells.Item(24, 19).Value = 0.001474614210218406
$ws.Cells.Item(24, 20).Value = 0.001063463620781133

$ws.Cells.Item(25, 5).Value = 2
$ws.Cells.Item(25, 7).Value = 38.730512
$ws.Cells.Item(25, 8).Value = 77.461024
$ws.Cells.Item(25, 9).Value = 0.00759933956842245
$ws.Cells.Item(25, 10).Value = 0.005079092294630384
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 13).Value = 17.78856566666667
$ws.Cells.Item(25, 14).Value = 53.365697
$ws.Cells.Item(25, 15).Value = 0.3777620057111716
$ws.Cells.Item(25, 16).Value = 0.4076169595435007
$ws.Cells.Item(25, 17).Value = 688.9602560156213
$ws.Cells.Item(25, 18).Value = 4133.761536093728
$ws.Cells.Item(25, 19).Value = 0.002870741757447534
$ws.Cells.Item(25, 20).Value = 0.002070324158378059

$ws.Cells.Item(26, 5).Value = 2
$ws.Cells.Item(26, 7).Value = 38.730512
$ws.Cells.Item(26, 8).Value = 77.461024
$ws.Cells.Item(26, 9).Value = 0.00759933956842245
$ws.Cells.Item(26, 10).Value = 0.005079092294630384
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 13).Value = 10.346848
$ws.Cells.Item(26, 14).Value = 20.693696
$ws.Cells.Item(26, 15).Value = 0.2197280054227695
$ws.Cells.Item(26, 16).Value = 0.1580622369691433
$ws.Cells.Item(26, 17).Value = 400.738720626176
$ws.Cells.Item(26, 18).Value = 1602.954882504704
$ws.Cells.Item(26, 19).Value = 0.001669787725899795
$ws.Cells.Item(26, 20).Value = 0.0008028126898620176
